# Automatic update of files.
# Rows 3 and 5 swap their species-observation data, and rows 6 and 7 swap
# their species-observation data (the underlying records were re-sorted /
# re-identified upstream; only the data columns move, the row's other
# metadata - location names, dates, reporters, etc. - stays where it is).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ArtRow($Row, $A, $B, $E, $F, $G, $H, $HasJ, $HasL, $M, $HasAF, $Q, $R) {
    $ws.Cells.Item($Row, 1).Value = $A   # A - Id
    $ws.Cells.Item($Row, 2).Value = $B   # B - Taxonsorteringsordning
    $ws.Cells.Item($Row, 5).Value = $E   # E - TaxonId
    $ws.Cells.Item($Row, 6).Value = $F   # F - Artnamn
    $ws.Cells.Item($Row, 7).Value = $G   # G - Vetenskapligt namn
    $ws.Cells.Item($Row, 8).Value = $H   # H - Auktor

    # Column I is always a blank placeholder cell on these rows; use it as a
    # template to (re)create other blank placeholder cells (J, L, AF) so the
    # underlying cell keeps existing even though it carries no value.
    $blankSource = $ws.Cells.Item($Row, 9)

    # J - Enhet (blank marker cell, present/absent depending on record)
    if ($HasJ) {
        $blankSource.Copy($ws.Cells.Item($Row, 10))
    } else {
        $ws.Cells.Item($Row, 10).Clear()
    }

    # L - Kon (blank marker cell, present/absent depending on record)
    if ($HasL) {
        $blankSource.Copy($ws.Cells.Item($Row, 12))
    } else {
        $ws.Cells.Item($Row, 12).Clear()
    }

    # M - Aktivitet
    if ($M -ne "") {
        $ws.Cells.Item($Row, 13).Value = $M
    } else {
        $ws.Cells.Item($Row, 13).Clear()
    }

    $ws.Cells.Item($Row, 17).Value = $Q  # Q - Ost
    $ws.Cells.Item($Row, 18).Value = $R  # R - Nord

    # AF - Bestamningsmetod (blank marker cell, present/absent depending on record)
    if ($HasAF) {
        $blankSource.Copy($ws.Cells.Item($Row, 32))
    } else {
        $ws.Cells.Item($Row, 32).Clear()
    }
}

# New row 3 <= old row 5's data (Garnlav / Alectoria sarmentosa)
Set-ArtRow 3 111741082 77515 6425 "Garnlav" "Alectoria sarmentosa" "(Ach.) Ach." $true $false "" $true 331468.5669229594 6627064.351006002

# New row 5 <= old row 3's data (Vedtrappmossa / Crossocalyx hellerianus)
Set-ArtRow 5 111741025 94134 53 "Vedtrappmossa" "Crossocalyx hellerianus" "(Nees ex Lindenb.) Meyl." $true $true "" $true 331437.2628167981 6627065.263253132

# New row 6 <= old row 7's data (Tretåig hackspett / Picoides tridactylus)
Set-ArtRow 6 111741120 56398 100109 "Tretåig hackspett" "Picoides tridactylus" "(Linnaeus, 1758)" $false $true "färska spår" $false 331468.5669229594 6627064.351006002

# New row 7 <= old row 6's data (Vedtrappmossa / Crossocalyx hellerianus)
Set-ArtRow 7 111741038 94134 53 "Vedtrappmossa" "Crossocalyx hellerianus" "(Nees ex Lindenb.) Meyl." $true $true "" $true 331443.3172632467 6627064.989183033

Write-Output "Row swap applied."
